$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 1.033196630495687
$ws.Range("Q3").Value = 1.897036734533635
$ws.Range("Q4").Value = 1.245846174113728
$ws.Range("Q5").Value = 1.784451737247451
$ws.Range("Q6").Value = 1.757553896613831
$ws.Range("Q7").Value = 2.267803543469689
$ws.Range("Q8").Value = 1.535174538115806
$ws.Range("Q9").Value = 1.439726994580897
$ws.Range("Q10").Value = 2.402440918079677
$ws.Range("Q11").Value = 2.402440918079677
$ws.Range("Q12").Value = 1.443841293612757
$ws.Range("Q13").Value = 1.211275628093407
$ws.Range("Q14").Value = 1.309435688241485
$ws.Range("Q15").Value = 1.384102258224339
$ws.Range("Q16").Value = 2.090989106311385
$ws.Range("Q17").Value = 1.264788290313442
$ws.Range("Q18").Value = 1.264788290313442
$ws.Range("Q19").Value = 1.79763062483054
$ws.Range("Q20").Value = 2.185032813561829
$ws.Range("Q21").Value = 1.251277561671179
$ws.Range("Q22").Value = 1.539342568313268
$ws.Range("Q23").Value = 1.194339638421116
$ws.Range("Q24").Value = 1.633643190573729
$ws.Range("Q25").Value = 1.342818542634564
$ws.Range("Q26").Value = 1.556506674803845
$ws.Range("Q27").Value = 1.939346032482382
$ws.Range("Q28").Value = 1.523639546448364
$ws.Range("Q29").Value = 1.533943525003013
$ws.Range("Q30").Value = 1.650922393499214
$ws.Range("Q31").Value = 1.650922393499214
$ws.Range("Q32").Value = 1.296420910492415
$ws.Range("Q33").Value = 1.296420910492415
$ws.Range("Q34").Value = 1.612698332942558
$ws.Range("Q35").Value = 1.450838758676033
$ws.Range("Q36").Value = 1.467604750719693
$ws.Range("Q37").Value = 1.689140930763304
$ws.Range("Q38").Value = 1.63699196957063
$ws.Range("Q39").Value = 1.451460227890586
$ws.Range("Q40").Value = 1.726090802103506
$ws.Range("Q41").Value = 1.726090802103506
$ws.Range("Q42").Value = 1.726090802103506
$ws.Range("Q43").Value = 1.704998939529603
$ws.Range("Q44").Value = 1.951846748730099
$ws.Range("Q45").Value = 1.977480319599974
$ws.Range("Q46").Value = 1.813585229043661
$ws.Range("Q47").Value = 1.636354764225714
$ws.Range("Q48").Value = 1.636354764225714
$ws.Range("Q49").Value = 1.679049255710149
$ws.Range("Q50").Value = 3.617164992432295
$ws.Range("Q51").Value = 2.188020963138853
$ws.Range("Q52").Value = 2.45663422781619
$ws.Range("Q53").Value = 2.45663422781619
$ws.Range("Q54").Value = 2.477646543958159
$ws.Range("Q55").Value = 2.723091978666176
$ws.Range("Q56").Value = 3.487009655055755
$ws.Range("Q57").Value = 3.165630025183241
$ws.Range("Q58").Value = 3.165630025183241
$ws.Range("Q59").Value = 3.198652056791207
$ws.Range("Q60").Value = 3.571290677517552
$ws.Range("Q61").Value = 3.824617860709115
$ws.Range("Q62").Value = 8.937750311733469
$ws.Range("Q63").Value = 8.937750311733469
